$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.919.81'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '3.033.90'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '593.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.89'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +7.19%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.027.44'
$ws.Range('E8').Value = '  +1.34%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.517'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.43'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.81%  '
$ws.Range('E11').Value = '  +2.21%  '
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.52'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.127'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('D16').Value = '3.536.09'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.09'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '62.915.23'
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('D19').Value = '3.032.46'
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '451.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.29'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.695'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  +1.74%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '83.08'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.31'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.83%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.03'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.57%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.14%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.52'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.70'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.21'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.55'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D35').Value = '0.0₃0865'
$ws.Range('E35').Value = '  +6.28%  '
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.92'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.99%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.16'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +10.89%  '
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.129'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +5.47%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '50.56'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.309'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +15.79%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.27'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +9.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '395.45'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D47').Value = '2.733.84'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.56'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('E50').Value = '  +3.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '24.36'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +4.19%  '
